$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week's price report was added for "Repollo" (Vega Monumental Concepción).
# This inserts two new data rows (Primera / Segunda quality) right above the
# existing row 84, shifting all the subsequent rows down by two (the sheet's
# used range grows from A1:R210 to A1:R212).
$ws.Rows("84:85").Insert()

# New row 84 - Crespo record, Primera
$ws.Range("A84").Value = 11
$ws.Range("B84").Value = "Vega Monumental Concepción"
$ws.Range("C84").Value = "Bíobío"
$ws.Range("D84").Value = 44477
$ws.Range("E84").Value = 8
$ws.Range("F84").Value = 100112006
$ws.Range("G84").Value = "Repollo"
$ws.Range("H84").Value = "Crespo record"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 700
$ws.Range("L84").Value = 800
$ws.Range("M84").Value = 750
$ws.Range("N84").Value = "$/unidad"
$ws.Range("O84").Value = "Región Metropolitana"
$ws.Range("P84").Value = 750
$ws.Range("Q84").Value = 1
$ws.Range("R84").Value = "Hortaliza"

# New row 85 - Crespo record, Segunda
$ws.Range("A85").Value = 11
$ws.Range("B85").Value = "Vega Monumental Concepción"
$ws.Range("C85").Value = "Bíobío"
$ws.Range("D85").Value = 44477
$ws.Range("E85").Value = 8
$ws.Range("F85").Value = 100112006
$ws.Range("G85").Value = "Repollo"
$ws.Range("H85").Value = "Crespo record"
$ws.Range("I85").Value = "Segunda"
$ws.Range("J85").Value = 500
$ws.Range("K85").Value = 600
$ws.Range("L85").Value = 600
$ws.Range("M85").Value = 600
$ws.Range("N85").Value = "$/unidad"
$ws.Range("O85").Value = "Región Metropolitana"
$ws.Range("P85").Value = 600
$ws.Range("Q85").Value = 1
$ws.Range("R85").Value = "Hortaliza"
